$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing first block (rows 4-9) ---
$ws.Range("B4").Value = 10000
$ws.Range("C4").Value = "Bit/s"
$ws.Range("B6").Value = 45
$ws.Range("B8").Formula = "=1000000/B4"

# --- New row 11: Real packet duration (experimental, first block) ---
$ws.Range("A10:D11").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4104) | Out-Null
$ws.Range("A4:C9").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null

$ws.Range("A11").Value = "Real packet duration"
$ws.Range("B11").Value = 62
$ws.Range("C11").Value = "ms"
$ws.Range("D11").Value = "Experimental"

$ws.Range("A14").Value = "Baudrate"
$ws.Range("B14").Value = 38400
$ws.Range("C14").Value = "Bit/s"

$ws.Range("A15").Value = "Packet length"
$ws.Range("B15").Value = 32
$ws.Range("C15").Value = "bytes"

$ws.Range("A16").Value = "Service info"
$ws.Range("B16").Value = 45
$ws.Range("C16").Value = "bytes"

$ws.Range("A17").Value = "Packet bit count"
$ws.Range("B17").Formula = "=8*(B15+B16)"
$ws.Range("C17").Value = "bits"

$ws.Range("A18").Value = "Bit length"
$ws.Range("B18").Formula = "=1000000/B14"
$ws.Range("C18").Value = "uS"

$ws.Range("A19").Value = "Packet length"
$ws.Range("B19").Formula = "=B17*B18/1000"
$ws.Range("C19").Value = "mS"
